# fix: unique command names in XLSX - prefix protocol name to each step
#
# For every "protocol" worksheet (i.e. every sheet other than the first five
# overview/summary sheets), prepend the sheet's own name + a space to the
# value of every cell in column A for the data rows (everything below the
# "Name" header row). This mirrors the author's change that made each
# command/step name unique across sheets by tagging it with its sheet name.

$wb = $excel.ActiveWorkbook

# Sheets that must NOT be touched (overview / summary sheets).
$excluded = @(
    "JockurworldJourney",
    "NRWaves",
    "PersonalJockurworld",
    "PositiveSpin",
    "ReEngagement"
)

foreach ($ws in $wb.Worksheets) {
    if ($excluded -contains $ws.Name) {
        continue
    }

    $used = $ws.UsedRange
    $firstRow = $used.Row
    $lastRow = $firstRow + $used.Rows.Count - 1

    # Row 1 is the header ("Name", "Text", "Note", "*Guidelines"); data starts
    # at row 2, so skip the header row if present.
    $startRow = [Math]::Max($firstRow, 2)

    $prefix = $ws.Name + " "

    for ($r = $startRow; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 1)
        $current = $cell.Value2
        if ($null -eq $current) {
            continue
        }
        $text = [string]$current
        if ($text.Length -eq 0) {
            continue
        }
        if ($text.StartsWith($prefix)) {
            # Already prefixed; avoid double-prefixing.
            continue
        }
        $cell.Value = $prefix + $text
    }
}
